$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.927.00'
$ws.Range('E2').Value = '  -0.04%  '
$ws.Range('D3').Value = '1.671.43'
$ws.Range('E3').Value = '  +1.21%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '214.62'
$ws.Range('D5').NumberFormat = 'General'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.11%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.518'
$ws.Range('D6').NumberFormat = 'General'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +1.08%  '
$ws.Range('E7').Value = '  -0.05%  '
$ws.Range('E8').Value = '  +0.03%  '
$ws.Range('E9').Value = '  +0.87%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '20.25'
$ws.Range('D10').NumberFormat = 'General'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.28%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0890'
$ws.Range('D11').NumberFormat = 'General'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +1.58%  '
$ws.Range('D12').Value = '1.907.48'
$ws.Range('E12').Value = '  +1.19%  '
$ws.Range('D13').Value = '1.643.12'
$ws.Range('E13').Value = '  -0.57%  '
$ws.Range('E14').Value = '  +0.15%  '
$ws.Range('E15').Value = '  +1.37%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '65.50'
$ws.Range('D16').NumberFormat = 'General'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.70%  '
$ws.Range('D17').Value = '26.929.82'
$ws.Range('E17').Value = '  -0.04%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '8.02'
$ws.Range('D18').NumberFormat = 'General'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +4.13%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '233.45'
$ws.Range('D19').NumberFormat = 'General'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.97%  '
$ws.Range('E20').Value = '  +0.17%  '
$ws.Range('E21').Value = '  -0.07%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.42'
$ws.Range('D22').NumberFormat = 'General'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.22%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '9.13'
$ws.Range('D23').NumberFormat = 'General'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -1.87%  '
$ws.Range('E24').Value = '  -1.93%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '146.33'
$ws.Range('D25').NumberFormat = 'General'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.68%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.11'
$ws.Range('D26').NumberFormat = 'General'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.23%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '15.92'
$ws.Range('D27').NumberFormat = 'General'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.95%  '
$ws.Range('E28').Value = '  -1.49%  '
$ws.Range('E29').Value = '  +0.05%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0497'
$ws.Range('D30').NumberFormat = 'General'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.20%  '
$ws.Range('E31').Value = '  +0.09%  '
$ws.Range('D33').Value = '1.456.26'
$ws.Range('E33').Value = '  -5.89%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.12'
$ws.Range('D34').NumberFormat = 'General'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +1.57%  '
$ws.Range('E35').Value = '  +2.12%  '
$ws.Range('E36').Value = '  +0.19%  '
$ws.Range('E37').Value = '  -0.70%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.899'
$ws.Range('D38').NumberFormat = 'General'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.84%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0170'
$ws.Range('D39').NumberFormat = 'General'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.70%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.05'
$ws.Range('D40').NumberFormat = 'General'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +12.95%  '
$ws.Range('E41').Value = '  -4.21%  '
$ws.Range('E42').Value = '  -0.04%  '
$ws.Range('E43').Value = '  +2.66%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '66.28'
$ws.Range('D44').NumberFormat = 'General'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.93%  '
$ws.Range('D45').Value = '1.811.48'
$ws.Range('E45').Value = '  +1.06%  '
$ws.Range('E46').Value = '  +0.74%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '90.74'
$ws.Range('D47').NumberFormat = 'General'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.40%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.53'
$ws.Range('D48').NumberFormat = 'General'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +1.09%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.101'
$ws.Range('D49').NumberFormat = 'General'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +2.40%  '
$ws.Range('E50').Value = '  +0.36%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '7.66'
$ws.Range('D51').NumberFormat = 'General'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.15%  '
